# Updates the "epl" column (column J) values on Sheet1 per resubmission
# recalculation (addl calcs, rerun all). Each cell is forced to Text format
# before assignment so the exact decimal-formatted string (including any
# trailing zeros, e.g. "0.010") is preserved instead of being reinterpreted
# as a number and losing precision/format in the General number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3 = "0.003"
    4 = "0.011"
    5 = "0.008"
    9 = "0.005"
    10 = "0.018"
    11 = "0.013"
    13 = "0.133"
    15 = "0.006"
    16 = "0.021"
    17 = "0.017"
    21 = "0.008"
    22 = "0.026"
    23 = "0.018"
    27 = "0.009"
    28 = "0.027"
    29 = "0.023"
    33 = "0.010"
    34 = "0.033"
    35 = "0.026"
    39 = "0.010"
    40 = "0.036"
    41 = "0.029"
    43 = "0.253"
    45 = "0.013"
    46 = "0.044"
    47 = "0.033"
    49 = "0.274"
    51 = "0.014"
    52 = "0.044"
    53 = "0.040"
    55 = "0.295"
    57 = "0.013"
    58 = "0.050"
    59 = "0.040"
    63 = "0.018"
    64 = "0.060"
    65 = "0.049"
    67 = "0.363"
    69 = "0.022"
    70 = "0.064"
    71 = "0.051"
    73 = "0.397"
    75 = "0.019"
    76 = "0.065"
    77 = "0.053"
    79 = "0.454"
    81 = "0.023"
    82 = "0.073"
    83 = "0.054"
    85 = "0.523"
    87 = "0.025"
    88 = "0.091"
    89 = "0.065"
    90 = "0.258"
    91 = "0.578"
    93 = "0.031"
    94 = "0.098"
    95 = "0.065"
    96 = "0.293"
    97 = "0.623"
    99 = "0.033"
    100 = "0.115"
    101 = "0.087"
    102 = "0.339"
    103 = "0.812"
    105 = "0.042"
    106 = "0.137"
    107 = "0.120"
    109 = "0.970"
    111 = "0.063"
    112 = "0.215"
    113 = "0.193"
    114 = "0.611"
    115 = "1.307"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 10)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
